# Actualización automática 2025-10-06 16:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D11").Value = 95.04000000000001
$ws1.Range("D23").Value = "1 de 21"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F11").Value = 95.04000000000001
$ws2.Range("F23").Value = 95.04000000000001

$ws2.Range("G3").Value = 0
$ws2.Range("G5").Value = 0
$ws2.Range("G6").Value = 0
$ws2.Range("G7").Value = 0
$ws2.Range("G8").Value = 0
$ws2.Range("G11").Value = 0
$ws2.Range("G13").Value = 0
$ws2.Range("G14").Value = 0
$ws2.Range("G15").Value = 0
$ws2.Range("G16").Value = 0
$ws2.Range("G17").Value = 0
$ws2.Range("G19").Value = 0
$ws2.Range("G21").Value = 0
$ws2.Range("G22").Value = 0
$ws2.Range("G23").Value = 0

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Column width in the raw OOXML is stored with a fixed 0.83 offset relative to
# the COM ColumnWidth property in this workbook's font metrics; 24.17 yields
# a stored width of 25 (matching the target diff).
$ws3.Columns.Item(6).ColumnWidth = 24.17

$ws3.Range("D3").Value = 95.04000000000001
$ws3.Range("E3").Value = 5409.57890386263
$ws3.Range("F3").Value = 0.01726550042062127

$ws3.Range("D14").Value = 95.04000000000001
$ws3.Range("E14").Value = 55329.70147880389
$ws3.Range("F14").Value = 0.001714757659922441
